# cater case for browser type in test data
$wb = $excel.ActiveWorkbook

# --- TestCases sheet: overall result flips from FAIL to PASS ---
$wsCases = $wb.Worksheets.Item("TestCases")
$wsCases.Range("D2").Value = "PASS"

# --- TestSteps sheet: fix the openbrowser(s) keyword + propagate PASS result ---
$wsSteps = $wb.Worksheets.Item("TestSteps")
$wsSteps.Range("E2").Value = "openbrowser"
$wsSteps.Range("G2").Value = "PASS"
$wsSteps.Range("G3").Value = "PASS"
$wsSteps.Range("G4").Value = "PASS"
$wsSteps.Range("G5").Value = "PASS"
$wsSteps.Range("G6").Value = "PASS"
$wsSteps.Range("G7").Value = "PASS"
$wsSteps.Range("G8").Value = "PASS"
[void]$wsSteps.Range("E21").Select()

# --- TestData sheet: capitalize the Chrome browser value used for this case ---
$wsData = $wb.Worksheets.Item("TestData")
$wsData.Range("D2").Value = "Chrome"
[void]$wsData.Range("D8").Select()

# restore the originally active sheet/tab
[void]$wsCases.Select()
